# Apply the "updated FOCUS SAX completion and added new Spine Box Casting
# Case print" edit to the "July 2018" log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")

# 1. FOCUS sax 100 Scale (row 10) finished printing on 23-07-2018 -> fill in
#    the "Date Completed" column.
$ws.Range("B10").Value = "23-07-2018"

# 2. New print job logged: "Spine Box Concept 2 Casting Case" on row 11.
$ws.Range("A11").Value = "23-07-2018"
$ws.Range("C11").Value = "Spine Box Concept 2 Casting Case"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "PCTPE"
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 20
$ws.Range("H11").Value = 0.2
$ws.Range("I11").Value = "NA"

# Leave the selection where the author left it when they saved the file.
[void]$ws.Range("H12").Select()
